# Add a "Save" column (H) to the s_vals sheet, matching the formatting
# of the existing header cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the neighboring header cell (G1) onto the new
# header cell (H1), then set its text.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data value for row 2.
$ws.Range("H2").Value = 1
